# Scheduled runner update: refresh market-board derived profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) for a handful of leve
# rows across the ALC/ARM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = ""
$ws.Range("N8").Value = ""

$ws.Range("H38").Value = 3403.5
$ws.Range("I38").Value = 302.5
$ws.Range("J38").Value = 6504.5
$ws.Range("K38").Value = 907.5
$ws.Range("L38").Value = 19513.5
$ws.Range("M38").Value = -535.5
$ws.Range("N38").Value = -20257.5

$ws.Range("H112").Value = 1567.2632
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1567.2632
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 4701.7896
$ws.Range("N112").Value = -6917.7896

$ws.Range("H116").Value = 3225.1333
$ws.Range("I116").Value = 2333
$ws.Range("J116").Value = 3671.2
$ws.Range("K116").Value = 2333
$ws.Range("L116").Value = 3671.2
$ws.Range("M116").Value = 1109
$ws.Range("N116").Value = -10555.2

$ws.Range("H132").Value = 2088.3333
$ws.Range("I132").Value = 1971
$ws.Range("J132").Value = 2323
$ws.Range("K132").Value = 5913
$ws.Range("L132").Value = 6969
$ws.Range("M132").Value = -3383
$ws.Range("N132").Value = -12029

$ws.Range("H137").Value = 2419.1353
$ws.Range("I137").Value = 1374.8518
$ws.Range("J137").Value = 5238.7
$ws.Range("K137").Value = 4124.555399999999
$ws.Range("L137").Value = 15716.1
$ws.Range("M137").Value = -1574.555399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3977.8708
$ws.Range("I32").Value = 3597.0688
$ws.Range("J32").Value = 9499.5
$ws.Range("K32").Value = 3597.0688
$ws.Range("L32").Value = 9499.5
$ws.Range("M32").Value = -3310.0688

$ws.Range("H61").Value = 1870.3334
$ws.Range("I61").Value = 1573.5385
$ws.Range("J61").Value = 3799.5
$ws.Range("K61").Value = 1573.5385
$ws.Range("L61").Value = 3799.5
$ws.Range("M61").Value = -1361.5385

$ws.Range("H74").Value = 15381185
$ws.Range("I74").Value = 24991300
$ws.Range("J74").Value = 4999.6
$ws.Range("K74").Value = 24991300
$ws.Range("L74").Value = 4999.6
$ws.Range("M74").Value = -24990426

$ws.Range("H77").Value = 15381185
$ws.Range("I77").Value = 24991300
$ws.Range("J77").Value = 4999.6
$ws.Range("K77").Value = 124956500
$ws.Range("L77").Value = 24998
$ws.Range("M77").Value = -124952132

$ws.Range("H110").Value = 2052.75
$ws.Range("I110").Value = 1626
$ws.Range("J110").Value = 3333
$ws.Range("K110").Value = 1626
$ws.Range("L110").Value = 3333
$ws.Range("M110").Value = 419

$ws.Range("H136").Value = 1870.3334
$ws.Range("I136").Value = 1573.5385
$ws.Range("J136").Value = 3799.5
$ws.Range("K136").Value = 4720.6155
$ws.Range("L136").Value = 11398.5
$ws.Range("M136").Value = -2170.6155

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4485.5
$ws.Range("I31").Value = 2357.25
$ws.Range("J31").Value = 12998.5
$ws.Range("K31").Value = 2357.25
$ws.Range("L31").Value = 12998.5
$ws.Range("M31").Value = -2062.25
$ws.Range("N31").Value = -13588.5

$ws.Range("H34").Value = 4485.5
$ws.Range("I34").Value = 2357.25
$ws.Range("J34").Value = 12998.5
$ws.Range("K34").Value = 2357.25
$ws.Range("L34").Value = 12998.5
$ws.Range("M34").Value = -2155.25
$ws.Range("N34").Value = -13402.5

$ws.Range("H58").Value = 2180.8
$ws.Range("I58").Value = 1439
$ws.Range("J58").Value = 2675.3333
$ws.Range("K58").Value = 1439
$ws.Range("L58").Value = 2675.3333
$ws.Range("M58").Value = -1236

$ws.Range("H134").Value = 1415.2142
$ws.Range("I134").Value = 1431.7693
$ws.Range("J134").Value = 1200
$ws.Range("K134").Value = 4295.3079
$ws.Range("L134").Value = 3600
$ws.Range("M134").Value = -1760.3079
$ws.Range("N134").Value = -8670

$ws.Range("H136").Value = 2180.8
$ws.Range("I136").Value = 1439
$ws.Range("J136").Value = 2675.3333
$ws.Range("K136").Value = 4317
$ws.Range("L136").Value = 8025.999899999999
$ws.Range("M136").Value = -1767

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 136.27272
$ws.Range("I7").Value = 88.77778000000001
$ws.Range("J7").Value = 350
$ws.Range("K7").Value = 266.33334
$ws.Range("L7").Value = 1050
$ws.Range("M7").Value = -154.33334

$ws.Range("H10").Value = 200
$ws.Range("I10").Value = 200
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 600
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -461

$ws.Range("H11").Value = 1154.25
$ws.Range("I11").Value = 344.9
$ws.Range("J11").Value = 5201
$ws.Range("K11").Value = 1034.7
$ws.Range("L11").Value = 15603
$ws.Range("M11").Value = -894.6999999999998
$ws.Range("N11").Value = -15883

$ws.Range("H92").Value = 683.3333
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 683.3333
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 2049.9999
$ws.Range("M92").Value = ""
$ws.Range("N92").Value = -4545.9999

$ws.Range("H94").Value = 2154.5
$ws.Range("I94").Value = 1410
$ws.Range("J94").Value = 2899
$ws.Range("K94").Value = 4230
$ws.Range("L94").Value = 8697
$ws.Range("M94").Value = -3554
$ws.Range("N94").Value = -10049

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 20008
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 20008
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 20008
$ws.Range("N29").Value = -20588

$ws.Range("H126").Value = 3597.5
$ws.Range("I126").Value = 3597.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 10792.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -8322.5
$ws.Range("N126").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2167500
$ws.Range("I2").Value = 2167500
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2167500
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -2167388

$ws.Range("H22").Value = 834.2
$ws.Range("I22").Value = 450
$ws.Range("J22").Value = 930.25
$ws.Range("K22").Value = 450
$ws.Range("L22").Value = 930.25
$ws.Range("M22").Value = -155
$ws.Range("N22").Value = -1520.25

$ws.Range("H27").Value = 834.2
$ws.Range("I27").Value = 450
$ws.Range("J27").Value = 930.25
$ws.Range("K27").Value = 450
$ws.Range("L27").Value = 930.25
$ws.Range("M27").Value = -343
$ws.Range("N27").Value = -1144.25

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = ""

$ws.Range("H55").Value = 298.85715
$ws.Range("I55").Value = 298.6
$ws.Range("J55").Value = 299.5
$ws.Range("K55").Value = 298.6
$ws.Range("L55").Value = 299.5
$ws.Range("M55").Value = -125.6
$ws.Range("N55").Value = -645.5

$ws.Range("H61").Value = 4916
$ws.Range("I61").Value = 4900.2
$ws.Range("J61").Value = 4995
$ws.Range("K61").Value = 4900.2
$ws.Range("L61").Value = 4995
$ws.Range("M61").Value = -4698.2
$ws.Range("N61").Value = -5399

$ws.Range("H113").Value = 4916
$ws.Range("I113").Value = 4900.2
$ws.Range("J113").Value = 4995
$ws.Range("K113").Value = 4900.2
$ws.Range("L113").Value = 4995
$ws.Range("M113").Value = -2730.2
$ws.Range("N113").Value = -9335

$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 3250
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 9750
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -14650

$ws.Range("H132").Value = 4649.1
$ws.Range("I132").Value = 4750
$ws.Range("J132").Value = 4623.875
$ws.Range("K132").Value = 14250
$ws.Range("L132").Value = 13871.625
$ws.Range("M132").Value = -11720
$ws.Range("N132").Value = -18931.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = ""

$ws.Range("H100").Value = 1091.125
$ws.Range("I100").Value = 1205
$ws.Range("J100").Value = 749.5
$ws.Range("K100").Value = 2410
$ws.Range("L100").Value = 1499
$ws.Range("M100").Value = -1869
$ws.Range("N100").Value = -2581

$ws.Range("H107").Value = 325.77777
$ws.Range("I107").Value = 366.7143
$ws.Range("J107").Value = 182.5
$ws.Range("K107").Value = 1100.1429
$ws.Range("L107").Value = 547.5
$ws.Range("M107").Value = 819.8571000000002

$ws.Range("H122").Value = 1961.5834
$ws.Range("I122").Value = 1961.5834
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5884.7502
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3434.7502

$ws.Range("H126").Value = 2019.6
$ws.Range("I126").Value = 1814.1428
$ws.Range("J126").Value = 2499
$ws.Range("K126").Value = 5442.428400000001
$ws.Range("L126").Value = 7497
$ws.Range("M126").Value = -2972.428400000001
$ws.Range("N126").Value = -12437

$ws.Range("H132").Value = 3054.8484
$ws.Range("I132").Value = 2373.8
$ws.Range("J132").Value = 4102.615
$ws.Range("K132").Value = 7121.400000000001
$ws.Range("L132").Value = 12307.845
$ws.Range("M132").Value = -4591.400000000001
